# 5.2.1.1b.xlsx — add a new "2022" data column (Q) mirroring the existing
# "2021" column (P): same formatting, new figures where the source has
# real numbers, same "…" placeholder where the source only has a dash.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone column P's formatting (styles/borders) into column Q for the
#    table body (rows 3-25) so every new cell picks up the right style
#    index, exactly like the original author would get from a
#    copy/paste-format of the previous year's column.
$ws.Range("P3:P25").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

# 2) Fill in the actual 2022 figures (only the rows that carry real data
#    in column P get real data in Q; the subtotal/header rows stay blank,
#    matching column P).
$ws.Range("Q4").Value = 2022

$ws.Range("Q5").Value = 8725
$ws.Range("Q7").Value = 8347
$ws.Range("Q8").Value = 378

# 3) Rows 10-25 in column P just show the "…" placeholder (no data broken
#    out by category yet) — mirror that into column Q.
$placeholder = [char]0x2026
$ws.Range("Q10").Value = $placeholder
$ws.Range("Q11").Value = $placeholder
$ws.Range("Q12").Value = $placeholder
$ws.Range("Q13").Value = $placeholder
$ws.Range("Q14").Value = $placeholder
$ws.Range("Q15").Value = $placeholder
$ws.Range("Q16").Value = $placeholder
$ws.Range("Q17").Value = $placeholder
$ws.Range("Q18").Value = $placeholder
$ws.Range("Q19").Value = $placeholder
$ws.Range("Q20").Value = $placeholder
$ws.Range("Q21").Value = $placeholder
$ws.Range("Q22").Value = $placeholder
$ws.Range("Q23").Value = $placeholder
$ws.Range("Q24").Value = $placeholder
$ws.Range("Q25").Value = $placeholder

# 4) Move the selection like the author's workbook ends up (one cell left
#    of where the new column's header used to sit before the edit).
$ws.Range("Q3").Select()
